$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2.5 GLB Anni Updates: insert 4 new dispatch entries before the old row 164
# ("Dream Interpretation" / Hecate), which shifts that row down to row 168.
$ws.Rows("164:167").Insert()

# Row 164: ID 1051001 - Expert at Work
$ws.Range("A164").Value = 1051001
$ws.Range("B164").Value = 'Purple'
$ws.Range("C164").Value = 'Expert at Work'
$ws.Range("D164").Value = '専門に合った仕事'
$ws.Range("E164").Value = '전문가'
$ws.Range("F164").Value = '专业对口'
$ws.Range("G164").Value = '專業對口'
$ws.Range("H164").Value = 'The FAC has detected an abnormal rise in M-value at a high school in Eastside. Preliminary investigation reveals that its students have been holding spooky tale-telling séances after school. A request has been submitted to the Bureau to dispatch a suitable Sinner to assist in uncovering the truth.'
$ws.Range("I164").Value = 'FACはニューシティのとある高校でM値が異常に上昇していることを検出した。初動調査の結果、その学校の生徒たちの間で、放課後の心霊話会が流行していることが判明した。そのため、真相を探るのに適切なコンビクトを派遣するよう管理局に依頼が来た。'
$ws.Range("J164").Value = 'FAC 관측 결과 신성의 모 학교 내에서 M 수치가 비정상적으로 상승했다고 한다. 조사를 진행해 보니 해당 학교의 학생들 사이에서 서로 신비한 이야기를 나누는 방과후 활동이 유행하고 있다고 한다. FAC는 특별히 관리국에 적합한 수감자를 파견해 진상조사를 도와 달라고 요청했다.'
$ws.Range("K164").Value = 'FAC监测到新城某高校内M值异常上升，初步调查得知，该校学生间正流行着举办通灵故事会的课后活动，特向管理局申请派一位合适的禁闭者协助探查真相。'
$ws.Range("L164").Value = 'FAC監測到新城某高校內M值異常上升，初步調查得知，該校學生間正流行著舉辦通靈故事會的課後活動，特向管理局申請派一位合適的禁閉者協助探查真相。'
$ws.Range("M164").Value = 'Yugu'
$ws.Range("N164").Value = '玉骨'
$ws.Range("O164").Value = '유구'
$ws.Range("P164").Value = '玉骨'
$ws.Range("Q164").Value = '玉骨'
$ws.Range("W164").Value = 'Infected Gel'
$ws.Range("X164").Value = '感染されたゲル'
$ws.Range("Y164").Value = '감염된 젤라틴'
$ws.Range("Z164").Value = '感染凝胶'
$ws.Range("AA164").Value = '感染凝膠'
$ws.Range("AB164").Value = '''1.0'
$ws.Range("AC164").Value = 'Organic Gel'
$ws.Range("AD164").Value = '原生ゲル'
$ws.Range("AE164").Value = '원시적 젤라틴'
$ws.Range("AF164").Value = '原生凝胶'
$ws.Range("AG164").Value = '原生凝膠'
$ws.Range("AH164").Value = '''1.0'

# Row 165: ID 1051002 - Heaven-Sent Karma
$ws.Range("A165").Value = 1051002
$ws.Range("B165").Value = 'Green'
$ws.Range("C165").Value = 'Heaven-Sent Karma'
$ws.Range("D165").Value = '天からの功徳'
$ws.Range("E165").Value = '하늘이 내린 공덕'
$ws.Range("F165").Value = '天降功德'
$ws.Range("G165").Value = '天降功德'
$ws.Range("H165").Value = 'The season is changing and flu outbreaks have emerged across regions. Recognizing this as a heaven-sent opportunity to accumulate good karma, a certain Sinner has eagerly requested permission to go out and provide free medical services.'
$ws.Range("I165").Value = '季節の変わり目に、多くの地域でインフルエンザが流行している。これが天からの功徳であると気付いたあるコンビクトは、外出して無料診察を行うことを強く希望した。'
$ws.Range("J165").Value = '환절기가 되면 여러 지역에서 유행성 감기가 발발하는데, 이를 하늘이 내린 공덕으로 여긴 한 수감자가 자선 의료 활동을 할 수 있도록 외출을 허가해달라고 강력히 요구했다.'
$ws.Range("K165").Value = '季节交替，多地爆发流感，意识到这是天降的功德，某位禁闭者强烈要求外出义诊。'
$ws.Range("L165").Value = '季節交替，多地爆發流感，意識到這是天降的功德，某位禁閉者強烈要求外出義診。'
$ws.Range("M165").Value = 'Wuhuanzi'
$ws.Range("N165").Value = '無患子'
$ws.Range("O165").Value = '무환자'
$ws.Range("P165").Value = '无患子'
$ws.Range("Q165").Value = '無患子'
$ws.Range("W165").Value = 'Organic Gel'
$ws.Range("X165").Value = '原生ゲル'
$ws.Range("Y165").Value = '원시적 젤라틴'
$ws.Range("Z165").Value = '原生凝胶'
$ws.Range("AA165").Value = '原生凝膠'
$ws.Range("AB165").Value = '''1.0'
$ws.Range("AC165").Value = 'Condense Gel'
$ws.Range("AD165").Value = '懸濁ゲル'
$ws.Range("AE165").Value = '현탁 젤라틴'
$ws.Range("AF165").Value = '悬浊凝胶'
$ws.Range("AG165").Value = '懸濁凝膠'
$ws.Range("AH165").Value = '''1.0'

# Row 166: ID 1051003 - Private Bodyguard
$ws.Range("A166").Value = 1051003
$ws.Range("B166").Value = 'Blue'
$ws.Range("C166").Value = 'Private Bodyguard'
$ws.Range("D166").Value = '同行する用心棒'
$ws.Range("E166").Value = '수행 경호원'
$ws.Range("F166").Value = '随行保镖'
$ws.Range("G166").Value = '隨行保鑣'
$ws.Range("H166").Value = 'SALVA has extended an invitation to discuss traditional Eastian medicine with a certain Sinner. As this Sinner is not a local, company of another Sinner of similar age who is familiar with the route from the Bureau to Syndicate is required.'
$ws.Range("I166").Value = 'あるコンビクトと東洲の現地医術について学術的な討論をしたいと、彼岸から招待状が届いた。しかし、そのコンビクトはディスの人間ではないため、管理局からシンジケートまでのルートに詳しく、年齢も近いコンビクトに同行してもらう必要がある。'
$ws.Range("J166").Value = '구원 병원에서 동방 대륙 본토 의술에 관해 특정 수감자와 함께 학술적인 연구를 진행하길 희망한다며 초대장을 보내왔다. 해당 수감자가 신디케이트 사람이 아니기에 관리국에서 신디케이트로 향하는 길에 익숙한 또래 나이의 수감자를 함께 파견해야 한다.'
$ws.Range("K166").Value = '彼岸发来邀请，希望能与某位禁闭者就东洲本土医术进行学术探讨，鉴于该禁闭者非本地人，需要一位熟悉从管理局到辛迪加路线且年龄相仿的禁闭者陪同上路。'
$ws.Range("L166").Value = '彼岸發來邀請，希望能與某位禁閉者就東洲本土醫術進行學術探討，鑒於該禁閉者非本地人，需要一位熟悉從管理局到辛迪加路線且年齡相仿的禁閉者陪同上路。'
$ws.Range("M166").Value = 'Wuhuanzi'
$ws.Range("N166").Value = '無患子'
$ws.Range("O166").Value = '무환자'
$ws.Range("P166").Value = '无患子'
$ws.Range("Q166").Value = '無患子'
$ws.Range("R166").Value = 'Dolly'
$ws.Range("S166").Value = 'ドリー'
$ws.Range("T166").Value = '돌리'
$ws.Range("U166").Value = '多莉'
$ws.Range("V166").Value = '多莉'
$ws.Range("W166").Value = 'Organic Gel'
$ws.Range("X166").Value = '原生ゲル'
$ws.Range("Y166").Value = '원시적 젤라틴'
$ws.Range("Z166").Value = '原生凝胶'
$ws.Range("AA166").Value = '原生凝膠'
$ws.Range("AB166").Value = '''1.0'
$ws.Range("AC166").Value = 'Organic Gel'
$ws.Range("AD166").Value = '原生ゲル'
$ws.Range("AE166").Value = '원시적 젤라틴'
$ws.Range("AF166").Value = '原生凝胶'
$ws.Range("AG166").Value = '原生凝膠'
$ws.Range("AH166").Value = '''1.0'

# Row 167: ID 1051004 - Parade Fairy
$ws.Range("A167").Value = 1051004
$ws.Range("B167").Value = 'Purple'
$ws.Range("C167").Value = 'Parade Fairy'
$ws.Range("D167").Value = 'フロートキャラクター'
$ws.Range("E167").Value = '퍼레이드 캐릭터'
$ws.Range("F167").Value = '花车仙子'
$ws.Range("G167").Value = '花車仙子'
$ws.Range("H167").Value = 'The Public Security Bureau has received an arson threat targeting an amusement park. The sender is suspected of possessing Mania weapons. A Sinner skilled in creating illusions is urgently required to go undercover as a fairy at the parade and assist the Public Security Bureau in controlling the situation while guiding visitors to safety.'
$ws.Range("I167").Value = '治安局に遊園地への放火予告の手紙が届いた。送り主は狂瞳武器を所持している疑いがあり、幻術に長けたコンビクトにフロートに乗ったキャラクターに変装してもらい、治安局と協力して現場を制御し、来場客を安全な場所へ誘導する必要がある。'
$ws.Range("J167").Value = '치안국에 놀이공원 방화를 예고하는 편지 한 통이 전달되었다. 편지를 보낸 이는 변이 무기를 소지한 것으로 추측된다. 이에 환술에 능한 수감자를 급히 파견해 퍼레이드 캐릭터로 분장시키고, 치안국을 도와 상황을 통제하며 관광객을 안전한 구역으로 피신시키려 한다.'
$ws.Range("K167").Value = '治安局收到一封游乐园纵火预告信，发信人疑似持有狂厄武器，急需一位精通幻术的禁闭者扮成花车仙子，协助治安局控制局面并将游客引导至安全地带。'
$ws.Range("L167").Value = '治安局收到一封遊樂園縱火預告信，發信人疑似持有狂厄武器，急需一位精通幻術的禁閉者扮成花車仙子，協助治安局控制局面並將遊客引導至安全地帶。'
$ws.Range("M167").Value = 'Yingying'
$ws.Range("N167").Value = 'インイン'
$ws.Range("O167").Value = '잉잉'
$ws.Range("P167").Value = '萦萦'
$ws.Range("Q167").Value = '縈縈'
$ws.Range("W167").Value = 'Infected Cyst'
$ws.Range("X167").Value = '感染された嚢胞'
$ws.Range("Y167").Value = '감염된 낭포'
$ws.Range("Z167").Value = '感染囊胞'
$ws.Range("AA167").Value = '感染囊胞'
$ws.Range("AB167").Value = '''1.0'
$ws.Range("AC167").Value = 'Organic Cyst'
$ws.Range("AD167").Value = '原生嚢胞'
$ws.Range("AE167").Value = '원시적 낭포'
$ws.Range("AF167").Value = '原生囊胞'
$ws.Range("AG167").Value = '原生囊胞'
$ws.Range("AH167").Value = '''1.0'

